# Update "Elast of Bldg Svc Demand wrt E Cost.xlsx" (EoBSDwEC) workbook.
#
# Changes:
#  1. "About" sheet: remove the trailing note "use US data" (row 24).
#  2. "EoBSDwEC" sheet: retitle the fuel-elasticity table header and add
#     four more fuel rows (kerosene, heavy/residual fuel oil,
#     LPG propane/butane, hydrogen) with the same Commercial/Residential
#     elasticity values used by the other non-electricity fuels.

$wb  = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. "About" sheet - drop the last row (B24 "use US data")
# ---------------------------------------------------------------------
$about = $wb.Worksheets.Item("About")
$about.Rows.Item(24).Delete()
[void]$about.Range("A1").Select()

# ---------------------------------------------------------------------
# 2. "EoBSDwEC" sheet - header + new fuel rows
# ---------------------------------------------------------------------
$data = $wb.Worksheets.Item("EoBSDwEC")

# Header cell A1: "Fuel" -> "Elasticity by Fuel (dimensionless)", bold + wrap
$data.Range("A1").Value2 = "Elasticity by Fuel (dimensionless)"
$data.Range("A1").Font.Bold = $true
$data.Range("A1").WrapText = $true
$data.Rows.Item(1).RowHeight = 30

# New fuel rows appended after "biomass" (row 7), each using the standard
# -0.15 / -0.15 / -0.25 elasticity values shared by the other fuels.
$newFuels = @(
    @{ Row = 8;  Name = "kerosene" },
    @{ Row = 9;  Name = "heavy or residual fuel oil" },
    @{ Row = 10; Name = "LPG propane or butane" },
    @{ Row = 11; Name = "hydrogen" }
)

foreach ($fuel in $newFuels) {
    $r = $fuel.Row
    $data.Cells.Item($r, 1).Value2 = $fuel.Name
    $data.Cells.Item($r, 2).Value2 = -0.15
    $data.Cells.Item($r, 3).Value2 = -0.15
    $data.Cells.Item($r, 4).Value2 = -0.25
}

[void]$data.Range("A1").Select()
